$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve column D as text so values like "154.20" or "0.515" are not
# coerced into floating point numbers (column has no special format in
# the source file; it stores human-formatted price strings).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.663.81'
$ws.Range("E2").Value = '  +4.11%  '
$ws.Range("D3").Value = '3.251.88'
$ws.Range("E3").Value = '  +7.20%  '
$ws.Range("D5").Value = '583.13'
$ws.Range("E5").Value = '  +5.12%  '
$ws.Range("D6").Value = '154.20'
$ws.Range("E6").Value = '  +9.37%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '3.241.80'
$ws.Range("E8").Value = '  +7.08%  '
$ws.Range("D9").Value = '0.515'
$ws.Range("E9").Value = '  +5.13%  '
$ws.Range("D10").Value = '7.12'
$ws.Range("E10").Value = '  +10.62%  '
$ws.Range("E11").Value = '  +6.17%  '
$ws.Range("D12").Value = '0.490'
$ws.Range("E12").Value = '  +5.00%  '
$ws.Range("D13").Value = '38.10'
$ws.Range("E13").Value = '  +4.38%  '
$ws.Range("E14").Value = '  +6.11%  '
$ws.Range("D15").Value = '3.771.68'
$ws.Range("E15").Value = '  +6.43%  '
$ws.Range("B16").Value = 'BitcoinCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D16").Value = '558.72'
$ws.Range("E16").Value = '  +13.61%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '66.671.62'
$ws.Range("E17").Value = '  +3.82%  '
$ws.Range("D18").Value = '3.251.49'
$ws.Range("E18").Value = '  +6.47%  '
$ws.Range("E19").Value = '  +3.28%  '
$ws.Range("D20").Value = '7.14'
$ws.Range("E20").Value = '  +6.57%  '
$ws.Range("D21").Value = '14.48'
$ws.Range("E21").Value = '  +5.31%  '
$ws.Range("D22").Value = '0.744'
$ws.Range("E22").Value = '  +7.85%  '
$ws.Range("D23").Value = '7.77'
$ws.Range("E23").Value = '  +8.44%  '
$ws.Range("D24").Value = '13.67'
$ws.Range("E24").Value = '  +7.42%  '
$ws.Range("D25").Value = '81.99'
$ws.Range("E25").Value = '  +3.74%  '
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '9.25'
$ws.Range("E27").Value = '  +18.36%  '
$ws.Range("D28").Value = '2.98'
$ws.Range("E28").Value = '  +8.01%  '
$ws.Range("D29").Value = '2.24'
$ws.Range("E29").Value = '  +6.51%  '
$ws.Range("D30").Value = '27.83'
$ws.Range("E30").Value = '  +6.83%  '
$ws.Range("D31").Value = '2.77'
$ws.Range("E31").Value = '  +5.61%  '
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("D33").Value = '1.19'
$ws.Range("E33").Value = '  +6.76%  '
$ws.Range("D34").Value = '562.22'
$ws.Range("E34").Value = '  +9.06%  '
$ws.Range("D35").Value = '5.73'
$ws.Range("E35").Value = '  +4.49%  '
$ws.Range("D36").Value = '6.39'
$ws.Range("E36").Value = '  +7.00%  '
$ws.Range("D37").Value = '0.0459'
$ws.Range("E37").Value = '  +13.53%  '
$ws.Range("D38").Value = '55.47'
$ws.Range("E38").Value = '  +5.46%  '
$ws.Range("D39").Value = '0.134'
$ws.Range("E39").Value = '  +9.54%  '
$ws.Range("D40").Value = '0.0863'
$ws.Range("E40").Value = '  +7.56%  '
$ws.Range("E41").Value = '  +14.32%  '
$ws.Range("D42").Value = '3.170.17'
$ws.Range("E42").Value = '  +8.65%  '
$ws.Range("D43").Value = '8.65'
$ws.Range("E43").Value = '  +3.08%  '
$ws.Range("D44").Value = '0.275'
$ws.Range("E44").Value = '  +11.80%  '
$ws.Range("D45").Value = '2.33'
$ws.Range("E45").Value = '  +9.85%  '
$ws.Range("D46").Value = '26.58'
$ws.Range("E46").Value = '  +5.38%  '
$ws.Range("D48").Value = '0.0₃0558'
$ws.Range("E48").Value = '  +4.06%  '
$ws.Range("D49").Value = '125.97'
$ws.Range("E49").Value = '  +5.25%  '
$ws.Range("E50").Value = '  +3.33%  '
$ws.Range("D51").Value = '2.24'
$ws.Range("E51").Value = '  +8.27%  '
